$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 1.306376666666667
$ws.Range("N2").Value = 3.91913
$ws.Range("O2").Value = 0.06159635513812315
$ws.Range("P2").Value = 0.07271399171915481
$ws.Range("Q2").Value = 0.6175738926466666
$ws.Range("R2").Value = 5.55816503382
$ws.Range("S2").Value = 0.002017680366056467
$ws.Range("T2").Value = 0.002470195229162786
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.4077922698431246
$ws.Range("P3").Value = 0.4813954277979023
$ws.Range("Q3").Value = 4.088583795478001
$ws.Range("R3").Value = 36.79725415930201
$ws.Range("S3").Value = 0.01335784324327383
$ws.Range("T3").Value = 0.01635367088193984
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.637617
$ws.Range("N4").Value = 1.912851
$ws.Range("O4").Value = 0.03006398091472189
$ws.Range("P4").Value = 0.03549028273468269
$ws.Range("Q4").Value = 0.301425785346
$ws.Range("R4").Value = 2.712832068114
$ws.Range("S4").Value = 0.0009847904779610474
$ws.Range("T4").Value = 0.001205654166689869
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 9.728125
$ws.Range("N5").Value = 19.45625
$ws.Range("O5").Value = 0.4586862714388558
$ws.Range("P5").Value = 0.3609835859963323
$ws.Range("Q5").Value = 4.59885435625
$ws.Range("R5").Value = 27.5931261375
$ws.Range("S5").Value = 0.01502495207689697
$ws.Range("T5").Value = 0.01226311347860328
$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 0.8878186666666666
$ws.Range("N6").Value = 2.663456
$ws.Range("O6").Value = 0.0418611226651744
$ws.Range("P6").Value = 0.0494167117519279
$ws.Range("Q6").Value = 0.4197056208426667
$ws.Range("R6").Value = 3.777350587584
$ws.Range("S6").Value = 0.001371223428938386
$ws.Range("T6").Value = 0.001678754290948501
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 1.306376666666667
$ws.Range("N7").Value = 3.91913
$ws.Range("O7").Value = 0.06159635513812315
$ws.Range("P7").Value = 0.07271399171915481
$ws.Range("Q7").Value = 15.50096039435556
$ws.Range("R7").Value = 139.5086435492
$ws.Range("S7").Value = 0.05064330570820307
$ws.Range("T7").Value = 0.06200132303113002
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.4077922698431246
$ws.Range("P8").Value = 0.4813954277979023
$ws.Range("S8").Value = 0.3352787440230455
$ws.Range("T8").Value = 0.4104733177059819
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.637617
$ws.Range("N9").Value = 1.912851
$ws.Range("O9").Value = 0.03006398091472189
$ws.Range("P9").Value = 0.03549028273468269
$ws.Range("Q9").Value = 7.56571677676
$ws.Range("R9").Value = 68.09145099084
$ws.Range("S9").Value = 0.02471801087671038
$ws.Range("T9").Value = 0.03026163785366142
$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 9.728125
$ws.Range("N10").Value = 19.45625
$ws.Range("O10").Value = 0.4586862714388558
$ws.Range("P10").Value = 0.3609835859963323
$ws.Range("Q10").Value = 115.4301697083333
$ws.Range("R10").Value = 692.5810182500001
$ws.Range("S10").Value = 0.3771227861866891
$ws.Range("T10").Value = 0.3078012827398998
$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 0.8878186666666666
$ws.Range("N11").Value = 2.663456
$ws.Range("O11").Value = 0.0418611226651744
$ws.Range("P11").Value = 0.0494167117519279
$ws.Range("Q11").Value = 10.53451300878222
$ws.Range("R11").Value = 94.81061707904
$ws.Range("S11").Value = 0.03441738764683686
$ws.Range("T11").Value = 0.04213634042126733
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 1.306376666666667
$ws.Range("N12").Value = 3.91913
$ws.Range("O12").Value = 0.06159635513812315
$ws.Range("P12").Value = 0.07271399171915481
$ws.Range("Q12").Value = 0.4950122465333333
$ws.Range("R12").Value = 4.4551102188
$ws.Range("S12").Value = 0.001617258279017376
$ws.Range("T12").Value = 0.001979968558132333
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.4077922698431246
$ws.Range("P13").Value = 0.4813954277979023
$ws.Range("Q13").Value = 3.27717715052
$ws.Range("R13").Value = 29.49459435468
$ws.Range("S13").Value = 0.01070689041655487
$ws.Range("T13").Value = 0.01310817613685518
$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.637617
$ws.Range("N14").Value = 1.912851
$ws.Range("O14").Value = 0.03006398091472189
$ws.Range("P14").Value = 0.03549028273468269
$ws.Range("Q14").Value = 0.24160583364
$ws.Range("R14").Value = 2.17445250276
$ws.Range("S14").Value = 0.0007893522583524067
$ws.Range("T14").Value = 0.0009663840792196205
$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 9.728125
$ws.Range("N15").Value = 19.45625
$ws.Range("O15").Value = 0.4586862714388558
$ws.Range("P15").Value = 0.3609835859963323
$ws.Range("Q15").Value = 3.686181125
$ws.Range("R15").Value = 22.11708675
$ws.Range("S15").Value = 0.01204315041519361
$ws.Range("T15").Value = 0.009829417054081444
$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 0.8878186666666666
$ws.Range("N16").Value = 2.663456
$ws.Range("O16").Value = 0.0418611226651744
$ws.Range("P16").Value = 0.0494167117519279
$ws.Range("Q16").Value = 0.3364122491733333
$ws.Range("R16").Value = 3.02771024256
$ws.Range("S16").Value = 0.001099095020271975
$ws.Range("T16").Value = 0.001345594337510854
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 1.306376666666667
$ws.Range("N17").Value = 3.91913
$ws.Range("O17").Value = 0.06159635513812315
$ws.Range("P17").Value = 0.07271399171915481
$ws.Range("Q17").Value = 2.02273745647
$ws.Range("R17").Value = 12.13642473882
$ws.Range("S17").Value = 0.006608500942479965
$ws.Range("T17").Value = 0.00539374745199357
$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.4077922698431246
$ws.Range("P18").Value = 0.4813954277979023
$ws.Range("Q18").Value = 13.391323184967
$ws.Range("R18").Value = 80.34793910980201
$ws.Range("S18").Value = 0.04375089392142314
$ws.Range("T18").Value = 0.03570874464043893
$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.637617
$ws.Range("N19").Value = 1.912851
$ws.Range("O19").Value = 0.03006398091472189
$ws.Range("P19").Value = 0.03549028273468269
$ws.Range("Q19").Value = 0.9872587452690001
$ws.Range("R19").Value = 5.923552471614
$ws.Range("S19").Value = 0.003225480562350252
$ws.Range("T19").Value = 0.002632583049629217
$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 9.728125
$ws.Range("N20").Value = 19.45625
$ws.Range("O20").Value = 0.4586862714388558
$ws.Range("P20").Value = 0.3609835859963323
$ws.Range("Q20").Value = 15.062610440625
$ws.Range("R20").Value = 60.25044176250001
$ws.Range("S20").Value = 0.04921116923735338
$ws.Range("T20").Value = 0.02677688641684504
$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 0.8878186666666666
$ws.Range("N21").Value = 2.663456
$ws.Range("O21").Value = 0.0418611226651744
$ws.Range("P21").Value = 0.0494167117519279
$ws.Range("Q21").Value = 1.374660247264
$ws.Range("R21").Value = 8.247961483584001
$ws.Range("S21").Value = 0.004491162958680605
$ws.Range("T21").Value = 0.003665611759114137
$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 1.306376666666667
$ws.Range("N22").Value = 3.91913
$ws.Range("O22").Value = 0.06159635513812315
$ws.Range("P22").Value = 0.07271399171915481
$ws.Range("Q22").Value = 0.2171981846
$ws.Range("R22").Value = 1.9547836614
$ws.Range("S22").Value = 0.0007096098423662753
$ws.Range("T22").Value = 0.0008687574487360965
$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.4077922698431246
$ws.Range("P23").Value = 0.4813954277979023
$ws.Range("Q23").Value = 1.43793801606
$ws.Range("R23").Value = 12.94144214454
$ws.Range("S23").Value = 0.004697898238827229
$ws.Range("T23").Value = 0.00575151843268643
$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.637617
$ws.Range("N24").Value = 1.912851
$ws.Range("O24").Value = 0.03006398091472189
$ws.Range("P24").Value = 0.03549028273468269
$ws.Range("Q24").Value = 0.10601020242
$ws.Range("R24").Value = 0.95409182178
$ws.Range("S24").Value = 0.0003463467393478073
$ws.Range("T24").Value = 0.0004240235854825665
$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 9.728125
$ws.Range("N25").Value = 19.45625
$ws.Range("O25").Value = 0.4586862714388558
$ws.Range("P25").Value = 0.3609835859963323
$ws.Range("Q25").Value = 1.6173980625
$ws.Range("R25").Value = 9.704388375000001
$ws.Range("S25").Value = 0.005284213522722713
$ws.Range("T25").Value = 0.004312886306902726
$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 0.8878186666666666
$ws.Range("N26").Value = 2.663456
$ws.Range("O26").Value = 0.0418611226651744
$ws.Range("P26").Value = 0.0494167117519279
$ws.Range("Q26").Value = 0.14760873152
$ws.Range("R26").Value = 1.32847858368
$ws.Range("S26").Value = 0.0004822536104465813
$ws.Range("T26").Value = 0.0005904109430870751
